$wb = $excel.ActiveWorkbook

# Use the "Germany" sheet as a structural template (11 rows, same column widths/styles)
# and copy it to the end of the workbook (after "Swiss"), then rename/update it for Portugal.
$template = $wb.Worksheets.Item("Germany")
$swiss = $wb.Worksheets.Item("Swiss")

$template.Copy([System.Reflection.Missing]::Value, $swiss)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Portugal"

$newSheet.Range("B2").Value = "Portugal Market"
$newSheet.Range("B4").Value = "NGC-3479/T2407"

$newSheet.Activate()
[void]$newSheet.Range("B4").Select()
